$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as text (e.g. with trailing zeros,
# thousands separators, etc.) - force Text format so Excel does not
# auto-convert the assigned strings into numbers and mangle them.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '29.549.84'
$ws.Range('E2').Value = '  -2.49%  '
$ws.Range('D3').Value = '2.003.75'
$ws.Range('E3').Value = '  -4.00%  '
$ws.Range('E4').Value = '  +0.99%  '
$ws.Range('D5').Value = '329.42'
$ws.Range('E5').Value = '  -3.85%  '
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('D7').Value = '0.5010'
$ws.Range('E7').Value = '  -4.18%  '
$ws.Range('D8').Value = '0.4224'
$ws.Range('E8').Value = '  -4.28%  '
$ws.Range('D9').Value = '54.23'
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('D10').Value = '0.09033'
$ws.Range('E10').Value = '  -3.03%  '
$ws.Range('D11').Value = '1.119'
$ws.Range('E11').Value = '  -4.07%  '
$ws.Range('D12').Value = '23.34'
$ws.Range('E12').Value = '  -5.78%  '
$ws.Range('D13').Value = '2.035.97'
$ws.Range('E13').Value = '  -2.27%  '
$ws.Range('D14').Value = '8.044'
$ws.Range('E14').Value = '  -6.22%  '
$ws.Range('D15').Value = '6.476'
$ws.Range('E15').Value = '  -5.99%  '
$ws.Range('D16').Value = '1.011'
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('D17').Value = '94.48'
$ws.Range('E17').Value = '  -6.50%  '
$ws.Range('D18').Value = '0.00001115'
$ws.Range('E18').Value = '  -3.75%  '
$ws.Range('D19').Value = '0.06680'
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('D20').Value = '19.70'
$ws.Range('E20').Value = '  -6.54%  '
$ws.Range('D21').Value = '1.012'
$ws.Range('E21').Value = '  +0.93%  '
$ws.Range('D22').Value = '5.975'
$ws.Range('E22').Value = '  -5.53%  '
$ws.Range('D23').Value = '29.593.57'
$ws.Range('E23').Value = '  -2.45%  '
$ws.Range('D24').Value = '12.01'
$ws.Range('E24').Value = '  -3.87%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').Value = '159.01'
$ws.Range('E26').Value = '  -2.09%  '
$ws.Range('D27').Value = '20.77'
$ws.Range('E27').Value = '  -4.56%  '
$ws.Range('D28').Value = '6.362'
$ws.Range('E28').Value = '  -4.58%  '
$ws.Range('D29').Value = '2.309'
$ws.Range('E29').Value = '  -8.05%  '
$ws.Range('D30').Value = '128.33'
$ws.Range('E30').Value = '  -3.43%  '
$ws.Range('D31').Value = '1.059'
$ws.Range('E31').Value = '  -6.52%  '
$ws.Range('D32').Value = '0.09974'
$ws.Range('E32').Value = '  -4.48%  '
$ws.Range('D33').Value = '1.570'
$ws.Range('E33').Value = '  -5.32%  '
$ws.Range('D34').Value = '5.842'
$ws.Range('E34').Value = '  -6.11%  '
$ws.Range('D35').Value = '3.796'
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('D36').Value = '0.02470'
$ws.Range('E36').Value = '  -5.81%  '
$ws.Range('D37').Value = '9.311'
$ws.Range('E37').Value = '  -8.00%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.06420'
$ws.Range('E38').Value = '  -6.00%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '1.305'
$ws.Range('E39').Value = '  -2.97%  '
$ws.Range('D40').Value = '0.6566'
$ws.Range('E40').Value = '  -5.65%  '
$ws.Range('D41').Value = '11.70'
$ws.Range('E41').Value = '  -6.27%  '
$ws.Range('D42').Value = '0.2054'
$ws.Range('E42').Value = '  -6.63%  '
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('D44').Value = '0.6369'
$ws.Range('E44').Value = '  -6.42%  '
$ws.Range('D45').Value = '13.43'
$ws.Range('E45').Value = '  -6.52%  '
$ws.Range('D46').Value = '2.204'
$ws.Range('E46').Value = '  -5.11%  '
$ws.Range('D47').Value = '1.305'
$ws.Range('E47').Value = '  -4.89%  '
$ws.Range('D48').Value = '3.513'
$ws.Range('E48').Value = '  -3.28%  '
$ws.Range('D49').Value = '0.00000000330'
$ws.Range('E49').Value = '  -3.98%  '
$ws.Range('D50').Value = '0.06990'
$ws.Range('E50').Value = '  -3.16%  '
$ws.Range('D51').Value = '1.130'
$ws.Range('E51').Value = '  -6.34%  '
